$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Project - Owner")
$ws2 = $wb.Worksheets.Item("Project - Invite - Member")

# --- Fix the stray test data value in sheet1 J21 (numberTotalContainers=9 -> =3) ---
$cur = $ws1.Range("J21").Value2
$new = $cur -replace 'numberTotalContainers=9$', 'numberTotalContainers=3'
$ws1.Range("J21").Value2 = $new

# --- Clear the STATUS (PASS/DEPFAIL) column L data rows on both sheets ---
$ws1.Range("L2:L48").ClearContents()
$ws2.Range("L2:L41").ClearContents()

# --- Update the view/selection state for sheet1 ("Project - Owner") ---
$ws1.Activate()
$ws1.Range("L2:L29").Select()

# --- Update the view/selection state for sheet2 ("Project - Invite - Member") ---
# Activating this sheet last makes it the active tab, matching the target state.
$ws2.Activate()
$ws2.Range("L2:L19").Select()
